$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report header text updates
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# Weekly crime statistics data refresh (rows 14-30)
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -20
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = 22.222222222222
$ws.Range("I14").Value = 13
$ws.Range("J14").Value = 13
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 62.5
$ws.Range("N14").Value = -76.363636363636
$ws.Range("C15").Value = 11
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 83.333333333333
$ws.Range("F15").Value = 27
$ws.Range("G15").Value = 35
$ws.Range("H15").Value = -22.857142857142
$ws.Range("I15").Value = 46
$ws.Range("J15").Value = 54
$ws.Range("K15").Value = -14.814814814814
$ws.Range("L15").Value = 43.75
$ws.Range("M15").Value = 39.393939393939
$ws.Range("N15").Value = -25.806451612903
$ws.Range("C16").Value = 64
$ws.Range("D16").Value = 82
$ws.Range("E16").Value = -21.951219512195
$ws.Range("F16").Value = 297
$ws.Range("G16").Value = 318
$ws.Range("H16").Value = -6.603773584905
$ws.Range("I16").Value = 505
$ws.Range("J16").Value = 488
$ws.Range("K16").Value = 3.483606557377
$ws.Range("L16").Value = 32.545931758530
$ws.Range("M16").Value = 3.271983640081
$ws.Range("N16").Value = -74.686716791979
$ws.Range("C17").Value = 120
$ws.Range("D17").Value = 144
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 482
$ws.Range("G17").Value = 462
$ws.Range("H17").Value = 4.329004329004
$ws.Range("I17").Value = 764
$ws.Range("J17").Value = 712
$ws.Range("K17").Value = 7.303370786516
$ws.Range("L17").Value = 26.699834162520
$ws.Range("M17").Value = 69.401330376940
$ws.Range("N17").Value = -10.011778563015
$ws.Range("D18").Value = 60
$ws.Range("E18").Value = 3.333333333333
$ws.Range("F18").Value = 225
$ws.Range("G18").Value = 222
$ws.Range("H18").Value = 1.351351351351
$ws.Range("I18").Value = 351
$ws.Range("J18").Value = 325
$ws.Range("K18").Value = 8
$ws.Range("L18").Value = 48.101265822784
$ws.Range("M18").Value = -12.25
$ws.Range("N18").Value = -84.167794316644
$ws.Range("C19").Value = 130
$ws.Range("D19").Value = 150
$ws.Range("E19").Value = -13.333333333333
$ws.Range("F19").Value = 524
$ws.Range("G19").Value = 565
$ws.Range("H19").Value = -7.256637168141
$ws.Range("I19").Value = 811
$ws.Range("J19").Value = 871
$ws.Range("K19").Value = -6.888633754305
$ws.Range("L19").Value = 29.967948717948
$ws.Range("M19").Value = 77.850877192982
$ws.Range("N19").Value = 2.012578616352
$ws.Range("C20").Value = 113
$ws.Range("D20").Value = 98
$ws.Range("E20").Value = 15.306122448979
$ws.Range("F20").Value = 408
$ws.Range("G20").Value = 355
$ws.Range("H20").Value = 14.929577464788
$ws.Range("I20").Value = 643
$ws.Range("J20").Value = 559
$ws.Range("K20").Value = 15.026833631484
$ws.Range("L20").Value = 188.340807174888
$ws.Range("M20").Value = 183.259911894273
$ws.Range("N20").Value = -65.149051490514
$ws.Range("C21").Value = 504
$ws.Range("D21").Value = 545
$ws.Range("E21").Value = -7.522935779816
$ws.Range("F21").Value = 1974
$ws.Range("G21").Value = 1966
$ws.Range("H21").Value = 0.406917599186
$ws.Range("I21").Value = 3133
$ws.Range("J21").Value = 3022
$ws.Range("K21").Value = 3.673064195896
$ws.Range("L21").Value = 48.272598201609
$ws.Range("M21").Value = 51.792635658914
$ws.Range("N21").Value = -59.925812228191
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 16
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = -45.238095238095
$ws.Range("L22").Value = -11.538461538461
$ws.Range("M22").Value = -28.125
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = 39
$ws.Range("E23").Value = -23.076923076923
$ws.Range("F23").Value = 121
$ws.Range("G23").Value = 117
$ws.Range("H23").Value = 3.418803418803
$ws.Range("I23").Value = 199
$ws.Range("J23").Value = 169
$ws.Range("K23").Value = 17.751479289940
$ws.Range("L23").Value = 59.2
$ws.Range("M23").Value = 80.909090909090
$ws.Range("C24").Value = 276
$ws.Range("D24").Value = 343
$ws.Range("E24").Value = -19.533527696793
$ws.Range("F24").Value = 1230
$ws.Range("G24").Value = 1245
$ws.Range("H24").Value = -1.204819277108
$ws.Range("I24").Value = 1846
$ws.Range("J24").Value = 1788
$ws.Range("K24").Value = 3.243847874720
$ws.Range("L24").Value = 28.910614525139
$ws.Range("M24").Value = 32.424677187948
$ws.Range("C25").Value = 178
$ws.Range("D25").Value = 189
$ws.Range("E25").Value = -5.820105820105
$ws.Range("F25").Value = 745
$ws.Range("G25").Value = 723
$ws.Range("H25").Value = 3.042876901798
$ws.Range("I25").Value = 1086
$ws.Range("J25").Value = 1032
$ws.Range("K25").Value = 5.232558139534
$ws.Range("L25").Value = 30.059880239521
$ws.Range("M25").Value = 0.835654596100
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 16.666666666666
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = -17.543859649122
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 87
$ws.Range("K26").Value = -12.643678160919
$ws.Range("L26").Value = 38.181818181818
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 22
$ws.Range("E27").Value = -31.818181818181
$ws.Range("F27").Value = 71
$ws.Range("G27").Value = 67
$ws.Range("H27").Value = 5.970149253731
$ws.Range("I27").Value = 125
$ws.Range("J27").Value = 89
$ws.Range("K27").Value = 40.449438202247
$ws.Range("L27").Value = 32.978723404255
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = -31.25
$ws.Range("F28").Value = 32
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = -21.951219512195
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 58
$ws.Range("K28").Value = -34.482758620689
$ws.Range("L28").Value = 8.571428571428
$ws.Range("M28").Value = -20.833333333333
$ws.Range("N28").Value = -69.354838709677
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 13
$ws.Range("E29").Value = -46.153846153846
$ws.Range("F29").Value = 23
$ws.Range("G29").Value = 36
$ws.Range("H29").Value = -36.111111111111
$ws.Range("I29").Value = 29
$ws.Range("J29").Value = 53
$ws.Range("K29").Value = -45.283018867924
$ws.Range("L29").Value = -9.375
$ws.Range("M29").Value = -32.558139534883
$ws.Range("N29").Value = -74.782608695652
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = -42.857142857142
$ws.Range("L30").Value = 33.333333333333

Write-Output "Applied weekly crime data update"
